$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 0.6160572926326269
$ws.Range("R2").Value = '[0.3686896860599518, 0.5712586641311646, 0.8399618864059448, 0.7188007235527039, 0.9159262180328369, 0.964373767375946, 0.8161059021949768]'
$ws.Range("T2").Value = '[0.3686896860599518, 0.964373767375946, 1.126319169998169, 1.2241731882095337, 1.2196450233459473, 1.4703971147537231, 1.557084560394287]'
$ws.Range("V2").Value = 0.5705868005752563
$ws.Range("Z2").Value = '2025-10-17T02:10:12.796669'

# Row 3
$ws.Range("F3").Value = 0.6160572926326269
$ws.Range("R3").Value = '[0.29108816385269165, 0.3007749617099762, 0.5837634801864624, 0.6349945664405823, 0.639224112033844, 0.5765426754951477, 0.6429972052574158]'
$ws.Range("T3").Value = '[0.639224112033844, 0.7289791107177734, 0.8478673100471497, 0.9712804555892944, 1.0848945379257202, 1.3480218648910522, 1.2196450233459473]'
$ws.Range("V3").Value = 0.537077784538269
$ws.Range("Z3").Value = '2025-10-17T02:10:12.796669'

# Row 4
$ws.Range("F4").Value = 0.6160572926326269
$ws.Range("R4").Value = '[1.5354762077331543, 1.5763767957687378, 1.5182515382766724, 1.6015644073486328, 1.814195990562439, 2.4762635231018066, 1.645581841468811]'
$ws.Range("T4").Value = '[2.4762635231018066, 2.9099297523498535, 2.5544352531433105, 3.1355085372924805, 2.66167950630188, 2.662851095199585, 2.9838945865631104]'
$ws.Range("V4").Value = 0.724717378616333
$ws.Range("Z4").Value = '2025-10-17T02:10:12.796669'

# Row 5
$ws.Range("F5").Value = 0.6160572926326269
$ws.Range("K5").Value = 0.4285714285714285
$ws.Range("R5").Value = '[0.4618493616580963, 0.46296268701553345, 0.7553806304931641, 0.8820913434028625, 0.8820913434028625, 1.2559460401535034, 1.8004459142684937]'
$ws.Range("T5").Value = '[0.7553806304931641, 2.353519916534424, 3.1355085372924805, 3.852335214614868, 2.665560245513916, 2.4172964096069336, 2.481671094894409]'
$ws.Range("V5").Value = 0.3056077361106873
$ws.Range("Z5").Value = '2025-10-17T02:10:12.796669'

# Row 6
$ws.Range("F6").Value = 0.6160572926326269
$ws.Range("R6").Value = '[0.9198511242866516, 0.9164473414421082, 1.0153549909591675, 0.964373767375946, 1.1352434158325195, 1.3090776205062866, 1.209598422050476]'
$ws.Range("T6").Value = '[0.964373767375946, 1.0153549909591675, 1.3090776205062866, 1.624017596244812, 1.7274311780929565, 1.9255344867706299, 1.6671695709228516]'
$ws.Range("V6").Value = 0.7060033082962036
$ws.Range("Z6").Value = '2025-10-17T02:10:12.797668'

# Row 7
$ws.Range("F7").Value = 0.6160572926326269
$ws.Range("R7").Value = '[1.0445706844329834, 0.7342439293861389, 1.3509747982025146, 0.9219695329666138, 0.9574594497680664, 1.027353048324585, 0.8839975595474243]'
$ws.Range("T7").Value = '[1.0150388479232788, 1.0873119831085205, 1.158658742904663, 1.3480218648910522, 1.4703971147537231, 1.486830234527588, 2.204462766647339]'
$ws.Range("V7").Value = 0.7479963302612305
$ws.Range("Z7").Value = '2025-10-17T02:10:12.797668'

# Row 8
$ws.Range("F8").Value = 0.6160572926326269
$ws.Range("R8").Value = '[0.2907591462135315, 0.28504446148872375, 0.3098526895046234, 0.4762324094772339, 0.6092221736907959, 0.693737268447876, 0.6997390985488892]'
$ws.Range("T8").Value = '[0.28504446148872375, 0.6743265986442566, 0.8478673100471497, 1.0150388479232788, 1.2241731882095337, 1.3018968105316162, 1.6027143001556396]'
$ws.Range("V8").Value = 0.3090955018997192
$ws.Range("Z8").Value = '2025-10-17T02:10:12.797668'

# Row 9
$ws.Range("F9").Value = 0.6160572926326269
$ws.Range("R9").Value = '[0.3686896860599518, 0.5519679188728333, 0.735078752040863, 0.7528611421585083, 1.0153549909591675, 0.8294775485992432, 1.004456877708435]'
$ws.Range("T9").Value = '[0.3686896860599518, 1.0153549909591675, 1.0848945379257202, 1.1840670108795166, 1.3018968105316162, 1.486830234527588, 1.4751609563827515]'
$ws.Range("V9").Value = 0.584369421005249
$ws.Range("Z9").Value = '2025-10-17T02:10:12.797668'

# Row 10
$ws.Range("F10").Value = 0.6160572926326269
$ws.Range("R10").Value = '[0.28504446148872375, 0.36131569743156433, 0.36825448274612427, 0.38687416911125183, 0.5760982632637024, 0.6288474202156067, 0.5670149326324463]'
$ws.Range("T10").Value = '[0.28504446148872375, 0.6349484324455261, 0.7289791107177734, 1.0873119831085205, 1.126319169998169, 1.1840670108795166, 1.4750926494598389]'
$ws.Range("V10").Value = 0.3809074759483337
$ws.Range("Z10").Value = '2025-10-17T02:10:12.799825'

# Row 11
$ws.Range("F11").Value = 0.6160572926326269
$ws.Range("R11").Value = '[0.39794325828552246, 0.3796297311782837, 0.5471352934837341, 0.5109684467315674, 0.4668952524662018, 0.5869107246398926, 0.639224112033844]'
$ws.Range("T11").Value = '[0.6349484324455261, 0.6743265986442566, 0.639224112033844, 1.1125237941741943, 1.158658742904663, 1.4751609563827515, 1.557084560394287]'
$ws.Range("V11").Value = 0.5016424655914307
$ws.Range("Z11").Value = '2025-10-17T02:10:12.800820'

# Row 12
$ws.Range("F12").Value = 0.6160572926326269
$ws.Range("R12").Value = '[1.022521734237671, 0.7418535351753235, 1.0917242765426636, 0.8169267773628235, 0.9712804555892944, 0.8513126373291016, 0.9974803328514099]'
$ws.Range("T12").Value = '[0.9712804555892944, 1.1125237941741943, 1.4750926494598389, 1.6027143001556396, 2.070220708847046, 1.913097858428955, 2.342046022415161]'
$ws.Range("V12").Value = 0.7366058826446533
$ws.Range("Z12").Value = '2025-10-17T02:10:12.800820'

# Row 13
$ws.Range("F13").Value = 0.6160572926326269
$ws.Range("R13").Value = '[1.961000919342041, 2.4762635231018066, 2.548816442489624, 3.9857804775238037, 2.8144288063049316, 2.781341075897217, 3.2413923740386963]'
$ws.Range("T13").Value = '[2.4762635231018066, 3.852335214614868, 3.940241575241089, 3.849271774291992, 3.7356345653533936, 4.5531134605407715, 4.46036958694458]'
$ws.Range("V13").Value = 0.9301201105117798
$ws.Range("Z13").Value = '2025-10-17T02:10:12.800820'

# Row 14
$ws.Range("F14").Value = 0.6160572926326269
$ws.Range("K14").Value = 0.4285714285714285
$ws.Range("R14").Value = '[0.4420248568058014, 0.4420248568058014, 0.7553806304931641, 0.8172926902770996, 1.0266729593276978, 1.2000261545181274, 1.2940304279327393]'
$ws.Range("T14").Value = '[0.7553806304931641, 2.204462766647339, 2.9099297523498535, 3.940241575241089, 2.5074117183685303, 2.7999069690704346, 2.5194039344787598]'
$ws.Range("V14").Value = 0.291774183511734
$ws.Range("Z14").Value = '2025-10-17T02:10:12.801823'

# Row 15
$ws.Range("F15").Value = 0.6160572926326269
$ws.Range("K15").Value = 0.5714285714285714
$ws.Range("R15").Value = '[0.8769850134849548, 0.9809556007385254, 0.998561441898346, 1.3458682298660278, 1.1135129928588867, 1.1277869939804077, 1.3090776205062866]'
$ws.Range("T15").Value = '[1.3090776205062866, 2.5544352531433105, 2.0157172679901123, 2.342046022415161, 1.981188416481018, 1.9998292922973633, 2.1122682094573975]'
$ws.Range("V15").Value = 0.6580042839050293
$ws.Range("Z15").Value = '2025-10-17T02:10:12.801823'

# Row 16
$ws.Range("F16").Value = 0.476247085334295
$ws.Range("R16").Value = '[3.548219919204712, 8.691213607788086, 5.758365631103516, 6.487487316131592, 9.90291976928711, 6.694366455078125, 5.858821868896484]'
$ws.Range("T16").Value = '[8.691213607788086, 7.731835842132568, 8.609739303588867, 8.581768989562988, 8.94638442993164, 9.074629783630371, 9.31938362121582]'
$ws.Range("V16").Value = 0.9216032028198242
$ws.Range("Z16").Value = '2025-10-17T02:10:13.019050'

# Row 17
$ws.Range("F17").Value = 0.476247085334295
$ws.Range("R17").Value = '[1.7238720655441284, 2.4434309005737305, 2.300889730453491, 2.6279022693634033, 2.213388204574585, 2.5505406856536865, 3.038969039916992]'
$ws.Range("T17").Value = '[1.7238720655441284, 2.300889730453491, 2.213388204574585, 2.4434309005737305, 2.6279022693634033, 3.038969039916992, 6.031808853149414]'
$ws.Range("V17").Value = 0.5427547693252563
$ws.Range("Z17").Value = '2025-10-17T02:10:13.019050'

# Row 18
$ws.Range("F18").Value = 0.476247085334295
$ws.Range("R18").Value = '[1.3550022840499878, 2.110875129699707, 2.9421937465667725, 2.3781371116638184, 2.6279022693634033, 2.9225618839263916, 3.3530654907226562]'
$ws.Range("T18").Value = '[2.6279022693634033, 2.7967379093170166, 2.9225618839263916, 3.6776535511016846, 3.725713014602661, 3.496691942214966, 5.9708662033081055]'
$ws.Range("V18").Value = 0.5590962171554565
$ws.Range("Z18").Value = '2025-10-17T02:10:13.020048'

# Row 19
$ws.Range("F19").Value = 0.476247085334295
$ws.Range("R19").Value = '[1.0273101329803467, 2.300889730453491, 2.022925615310669, 2.9225618839263916, 3.8935210704803467, 3.228111982345581, 4.250108242034912]'
$ws.Range("T19").Value = '[1.0273101329803467, 2.022925615310669, 2.300889730453491, 2.9225618839263916, 3.228111982345581, 3.639315605163574, 5.8619890213012695]'
$ws.Range("V19").Value = 0.5870458483695984
$ws.Range("Z19").Value = '2025-10-17T02:10:13.020048'

# Row 20
$ws.Range("F20").Value = 0.476247085334295
$ws.Range("R20").Value = '[1.2703914642333984, 2.022925615310669, 2.2484099864959717, 2.213388204574585, 2.7271506786346436, 2.191220998764038, 2.463343381881714]'
$ws.Range("T20").Value = '[2.022925615310669, 2.213388204574585, 2.2484099864959717, 2.7271506786346436, 2.7967379093170166, 3.532566547393799, 5.384163856506348]'
$ws.Range("V20").Value = 0.5372540950775146
$ws.Range("Z20").Value = '2025-10-17T02:10:13.020048'

# Row 21
$ws.Range("F21").Value = 0.476247085334295
$ws.Range("R21").Value = '[8.691213607788086, 9.207416534423828, 8.894898414611816, 13.830842971801758, 10.461819648742676, 11.372812271118164, 9.74061107635498]'
$ws.Range("T21").Value = '[8.691213607788086, 9.828558921813965, 9.933645248413086, 10.808511734008789, 11.14076042175293, 11.009100914001465, 10.901507377624512]'
$ws.Range("V21").Value = 1.256728887557983
$ws.Range("Z21").Value = '2025-10-17T02:10:13.020048'

# Row 22
$ws.Range("F22").Value = 0.476247085334295
$ws.Range("R22").Value = '[1.0273101329803467, 2.4434309005737305, 3.7905890941619873, 3.322721242904663, 2.7271506786346436, 4.417076587677002, 3.6776535511016846]'
$ws.Range("T22").Value = '[1.0273101329803467, 2.4434309005737305, 2.7271506786346436, 3.322721242904663, 3.6776535511016846, 3.9307637214660645, 6.2379069328308105]'
$ws.Range("V22").Value = 0.6351715326309204
$ws.Range("Z22").Value = '2025-10-17T02:10:13.020048'

# Row 23
$ws.Range("F23").Value = 0.476247085334295
$ws.Range("R23").Value = '[3.9940896034240723, 3.612931966781616, 3.9857804775238037, 3.513723611831665, 3.133359909057617, 3.7720181941986084, 8.738912582397461]'
$ws.Range("T23").Value = '[3.612931966781616, 5.384163856506348, 5.8619890213012695, 5.9708662033081055, 6.016942501068115, 6.031808853149414, 6.2379069328308105]'
$ws.Range("V23").Value = 0.887856125831604
$ws.Range("Z23").Value = '2025-10-17T02:10:13.020048'

# Row 24
$ws.Range("F24").Value = 0.476247085334295
$ws.Range("R24").Value = '[3.3579466342926025, 3.038969039916992, 3.612931966781616, 3.0922813415527344, 2.8238894939422607, 3.496691942214966, 3.9307637214660645]'
$ws.Range("T24").Value = '[3.038969039916992, 3.496691942214966, 3.612931966781616, 3.639315605163574, 3.699251413345337, 3.532566547393799, 3.9307637214660645]'
$ws.Range("V24").Value = 0.9351437091827393
$ws.Range("Z24").Value = '2025-10-17T02:10:13.021048'

# Row 25
$ws.Range("F25").Value = 0.476247085334295
$ws.Range("R25").Value = '[1.7238720655441284, 2.2484099864959717, 3.682757616043091, 3.322721242904663, 2.8477210998535156, 3.228111982345581, 3.136768102645874]'
$ws.Range("T25").Value = '[1.7238720655441284, 2.2484099864959717, 3.322721242904663, 3.228111982345581, 3.725713014602661, 3.699251413345337, 6.016942501068115]'
$ws.Range("V25").Value = 0.6654731035232544
$ws.Range("Z25").Value = '2025-10-17T02:10:13.021048'

# Row 26
$ws.Range("F26").Value = 0.5844655394036213
$ws.Range("K26").Value = 0.7142857142857143
$ws.Range("R26").Value = '[3.682757616043091, 5.021600246429443, 4.706303119659424, 9.207416534423828, 5.405718803405762, 8.388983726501465, 7.5876994132995605]'
$ws.Range("T26").Value = '[6.087997913360596, 7.964256763458252, 8.471141815185547, 7.782583713531494, 5.954493045806885, 9.602742195129395, 6.843430519104004]'
$ws.Range("V26").Value = 1.090826153755188
$ws.Range("Z26").Value = '2025-10-17T02:10:13.177216'

# Row 27
$ws.Range("F27").Value = 0.5844655394036213
$ws.Range("K27").Value = 0.7142857142857143
$ws.Range("R27").Value = '[0.8132038116455078, 1.586868405342102, 1.5875941514968872, 1.4857956171035767, 1.6002230644226074, 1.9224257469177246, 1.6309336423873901]'
$ws.Range("T27").Value = '[3.809053421020508, 3.489682912826538, 6.087997913360596, 5.247391700744629, 5.183218955993652, 4.85034704208374, 7.00767183303833]'
$ws.Range("V27").Value = 0.587891697883606
$ws.Range("Z27").Value = '2025-10-17T02:10:13.177216'

# Row 28
$ws.Range("F28").Value = 0.5844655394036213
$ws.Range("R28").Value = '[0.6938520073890686, 0.7095147371292114, 1.2714117765426636, 1.5296647548675537, 1.3206110000610352, 1.1889803409576416, 1.4472988843917847]'
$ws.Range("T28").Value = '[0.6938520073890686, 0.7095147371292114, 1.2714117765426636, 1.5296647548675537, 1.3206110000610352, 1.1889803409576416, 1.4472988843917847]'
$ws.Range("V28").Value = 0.6609089970588684
$ws.Range("Z28").Value = '2025-10-17T02:10:13.178215'

# Row 29
$ws.Range("F29").Value = 0.5844655394036213
$ws.Range("R29").Value = '[1.1962192058563232, 1.4052584171295166, 1.5114425420761108, 1.5090450048446655, 1.5793581008911133, 1.7265819311141968, 2.7678287029266357]'
$ws.Range("T29").Value = '[1.1962192058563232, 1.4052584171295166, 1.5114425420761108, 1.5090450048446655, 1.5793581008911133, 1.7265819311141968, 2.7678287029266357]'
$ws.Range("V29").Value = 0.6468383073806763
$ws.Range("Z29").Value = '2025-10-17T02:10:13.178215'

# Row 30
$ws.Range("F30").Value = 0.5844655394036213
$ws.Range("R30").Value = '[0.5247544646263123, 1.0627883672714233, 1.0364985466003418, 1.026412844657898, 1.042699933052063, 1.2714117765426636, 1.5825912952423096]'
$ws.Range("T30").Value = '[0.5247544646263123, 1.0627883672714233, 1.0364985466003418, 1.026412844657898, 1.042699933052063, 1.2714117765426636, 1.5825912952423096]'
$ws.Range("V30").Value = 0.6342445611953735
$ws.Range("Z30").Value = '2025-10-17T02:10:13.178215'

# Row 31
$ws.Range("F31").Value = 0.5844655394036213
$ws.Range("R31").Value = '[0.5955532789230347, 0.7853057980537415, 1.0364985466003418, 1.1688058376312256, 1.3206110000610352, 1.215668797492981, 1.6141602993011475]'
$ws.Range("T31").Value = '[0.5955532789230347, 0.7853057980537415, 1.0364985466003418, 1.1688058376312256, 1.3206110000610352, 1.215668797492981, 1.6141602993011475]'
$ws.Range("V31").Value = 0.5639538764953613
$ws.Range("Z31").Value = '2025-10-17T02:10:13.178215'

# Row 32
$ws.Range("F32").Value = 0.5844655394036213
$ws.Range("R32").Value = '[4.054708003997803, 2.8743762969970703, 3.549809217453003, 4.931705474853516, 3.940061569213867, 3.491143226623535, 4.050656318664551]'
$ws.Range("T32").Value = '[3.809053421020508, 4.054708003997803, 4.931705474853516, 4.90034294128418, 8.471141815185547, 5.014385223388672, 5.741334438323975]'
$ws.Range("V32").Value = 0.9614191055297852
$ws.Range("Z32").Value = '2025-10-17T02:10:13.178215'

# Row 33
$ws.Range("F33").Value = 0.5844655394036213
$ws.Range("R33").Value = '[0.5768178105354309, 0.5955532789230347, 1.0627883672714233, 1.151566743850708, 1.4561500549316406, 1.6920398473739624, 1.856772780418396]'
$ws.Range("T33").Value = '[0.5768178105354309, 0.5955532789230347, 1.0627883672714233, 1.151566743850708, 1.4561500549316406, 1.6920398473739624, 1.856772780418396]'
$ws.Range("V33").Value = 0.481997549533844
$ws.Range("Z33").Value = '2025-10-17T02:10:13.179213'

# Row 34
$ws.Range("F34").Value = 0.5844655394036213
$ws.Range("R34").Value = '[0.963847815990448, 1.042699933052063, 1.151566743850708, 1.2844855785369873, 1.1688058376312256, 1.1889803409576416, 1.2921749353408813]'
$ws.Range("T34").Value = '[0.963847815990448, 1.042699933052063, 1.151566743850708, 1.2844855785369873, 1.1688058376312256, 1.1889803409576416, 1.2921749353408813]'
$ws.Range("V34").Value = 0.7362955212593079
$ws.Range("Z34").Value = '2025-10-17T02:10:13.179213'

# Row 35
$ws.Range("F35").Value = 0.5844655394036213
$ws.Range("R35").Value = '[0.6938520073890686, 1.2321480512619019, 1.2844855785369873, 1.5825912952423096, 1.6141602993011475, 1.9720324277877808, 1.7865649461746216]'
$ws.Range("T35").Value = '[0.6938520073890686, 1.2321480512619019, 1.2844855785369873, 1.5825912952423096, 1.6141602993011475, 1.9720324277877808, 1.7865649461746216]'
$ws.Range("V35").Value = 0.7882651090621948
$ws.Range("Z35").Value = '2025-10-17T02:10:13.179213'

# Row 36
$ws.Range("F36").Value = 0.5844655394036213
$ws.Range("K36").Value = 0.4285714285714285
$ws.Range("R36").Value = '[2.265591859817505, 2.0094025135040283, 2.374429702758789, 2.0699026584625244, 2.1498682498931885, 2.195941925048828, 2.633648633956909]'
$ws.Range("T36").Value = '[2.265591859817505, 3.489682912826538, 2.0094025135040283, 2.374429702758789, 2.0699026584625244, 2.195941925048828, 2.442943572998047]'
$ws.Range("V36").Value = 1.053199052810669
$ws.Range("Z36").Value = '2025-10-17T02:10:13.179213'

# Row 37
$ws.Range("F37").Value = 0.5844655394036213
$ws.Range("R37").Value = '[0.5542061924934387, 0.9788116216659546, 1.1962192058563232, 1.1553232669830322, 1.1788504123687744, 1.1968179941177368, 2.302386522293091]'
$ws.Range("T37").Value = '[0.5542061924934387, 0.9788116216659546, 1.1962192058563232, 1.1553232669830322, 1.1788504123687744, 1.1968179941177368, 2.302386522293091]'
$ws.Range("V37").Value = 0.4933568239212036
$ws.Range("Z37").Value = '2025-10-17T02:10:13.180214'

# Row 38
$ws.Range("F38").Value = 0.5844655394036213
$ws.Range("K38").Value = 0.4285714285714285
$ws.Range("R38").Value = '[4.337198257446289, 4.054708003997803, 6.487487316131592, 5.44020938873291, 8.82075309753418, 5.109254360198975, 4.3913702964782715]'
$ws.Range("T38").Value = '[4.054708003997803, 4.337198257446289, 5.247391700744629, 5.44020938873291, 4.3913702964782715, 4.788558006286621, 4.678295612335205]'
$ws.Range("V38").Value = 1.074542284011841
$ws.Range("Z38").Value = '2025-10-17T02:10:13.180403'

# Row 39
$ws.Range("F39").Value = 0.5844655394036213
$ws.Range("R39").Value = '[4.087904453277588, 5.758365631103516, 3.762176036834717, 3.718254804611206, 5.44020938873291, 3.840088129043579, 3.842573881149292]'
$ws.Range("T39").Value = '[4.087904453277588, 3.762176036834717, 5.44020938873291, 3.718254804611206, 3.840088129043579, 3.842573881149292, 3.7610156536102295]'
$ws.Range("V39").Value = 1.150546550750732
$ws.Range("Z39").Value = '2025-10-17T02:10:13.180403'

# Row 40
$ws.Range("F40").Value = 0.5844655394036213
$ws.Range("R40").Value = '[0.5247544646263123, 0.5768178105354309, 0.7853057980537415, 0.963847815990448, 1.1512864828109741, 1.4472988843917847, 1.6623544692993164]'
$ws.Range("T40").Value = '[0.5247544646263123, 0.5768178105354309, 0.7853057980537415, 0.963847815990448, 1.1512864828109741, 1.4472988843917847, 1.6623544692993164]'
$ws.Range("V40").Value = 0.4038631319999695
$ws.Range("Z40").Value = '2025-10-17T02:10:13.180939'

# Row 41
$ws.Range("F41").Value = 0.5844655394036213
$ws.Range("K41").Value = 0.5714285714285714
$ws.Range("R41").Value = '[3.548219919204712, 4.337198257446289, 4.087904453277588, 4.931705474853516, 8.894898414611816, 4.53907585144043, 4.718254566192627]'
$ws.Range("T41").Value = '[4.337198257446289, 4.087904453277588, 4.931705474853516, 5.183218955993652, 4.53907585144043, 4.718254566192627, 4.522115707397461]'
$ws.Range("V41").Value = 0.9351364374160767
$ws.Range("Z41").Value = '2025-10-17T02:10:13.180939'

# Row 42
$ws.Range("F42").Value = 0.5844655394036213
$ws.Range("R42").Value = '[0.4555348753929138, 0.5437226295471191, 0.8630491495132446, 0.9788116216659546, 0.9866169691085815, 1.5090450048446655, 1.987839698791504]'
$ws.Range("T42").Value = '[0.4555348753929138, 0.5437226295471191, 0.8630491495132446, 0.9788116216659546, 0.9866169691085815, 1.5090450048446655, 1.987839698791504]'
$ws.Range("V42").Value = 0.3561474978923798
$ws.Range("Z42").Value = '2025-10-17T02:10:13.180939'

# Row 43
$ws.Range("F43").Value = 0.5844655394036213
$ws.Range("R43").Value = '[0.8630491495132446, 0.9548653364181519, 1.1553232669830322, 1.0304895639419556, 1.2386415004730225, 1.5793581008911133, 2.4712913036346436]'
$ws.Range("T43").Value = '[0.8630491495132446, 0.9548653364181519, 1.1553232669830322, 1.0304895639419556, 1.2386415004730225, 1.5793581008911133, 2.4712913036346436]'
$ws.Range("V43").Value = 0.4929726719856262
$ws.Range("Z43").Value = '2025-10-17T02:10:13.180939'

# Row 44
$ws.Range("F44").Value = 0.5844655394036213
$ws.Range("R44").Value = '[1.225460410118103, 1.5296647548675537, 1.9720324277877808, 1.7913583517074585, 2.302386522293091, 2.007189989089966, 2.4712913036346436]'
$ws.Range("T44").Value = '[1.225460410118103, 1.5296647548675537, 1.9720324277877808, 1.7913583517074585, 2.302386522293091, 2.007189989089966, 2.4712913036346436]'
$ws.Range("V44").Value = 0.8767750263214111
$ws.Range("Z44").Value = '2025-10-17T02:10:13.181927'

# Row 45
$ws.Range("F45").Value = 0.5844655394036213
$ws.Range("R45").Value = '[0.5542061924934387, 0.9866169691085815, 1.1374363899230957, 1.0372637510299683, 1.2386415004730225, 1.4052584171295166, 2.265591859817505]'
$ws.Range("T45").Value = '[0.5542061924934387, 0.9866169691085815, 1.1374363899230957, 1.0372637510299683, 1.2386415004730225, 1.4052584171295166, 2.265591859817505]'
$ws.Range("V45").Value = 0.5129446983337402
$ws.Range("Z45").Value = '2025-10-17T02:10:13.181927'

# Row 46
$ws.Range("F46").Value = 0.5844655394036213
$ws.Range("R46").Value = '[0.49568384885787964, 0.5437226295471191, 1.0304895639419556, 1.0372637510299683, 1.1968179941177368, 1.7265819311141968, 1.6623544692993164]'
$ws.Range("T46").Value = '[0.49568384885787964, 0.5437226295471191, 1.0304895639419556, 1.0372637510299683, 1.1968179941177368, 1.7265819311141968, 1.6623544692993164]'
$ws.Range("V46").Value = 0.4150242209434509
$ws.Range("Z46").Value = '2025-10-17T02:10:13.181927'

# Row 47
$ws.Range("F47").Value = 0.5844655394036213
$ws.Range("R47").Value = '[0.4555348753929138, 0.49568384885787964, 0.9548653364181519, 1.1374363899230957, 1.1788504123687744, 1.5114425420761108, 2.0850367546081543]'
$ws.Range("T47").Value = '[0.4555348753929138, 0.49568384885787964, 0.9548653364181519, 1.1374363899230957, 1.1788504123687744, 1.5114425420761108, 2.442943572998047]'
$ws.Range("V47").Value = 0.3568399250507355
$ws.Range("Z47").Value = '2025-10-17T02:10:13.181927'

# Row 48
$ws.Range("F48").Value = 0.5844655394036213
$ws.Range("R48").Value = '[0.7095147371292114, 1.225460410118103, 1.026412844657898, 1.2321480512619019, 1.1512864828109741, 1.215668797492981, 1.4561500549316406]'
$ws.Range("T48").Value = '[0.7095147371292114, 1.225460410118103, 1.026412844657898, 1.2321480512619019, 1.1512864828109741, 1.215668797492981, 1.4561500549316406]'
$ws.Range("V48").Value = 0.7167636752128601
$ws.Range("Z48").Value = '2025-10-17T02:10:13.181927'

Write-Host "Applied all cell updates."